# Updated cryptos list on Thu Nov 30 13:38:53 UTC 2023 with GitHub Actions
# Refresh crypto price/volume snapshot values in the ranking table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.905.79"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.046.21"
$ws.Range("E3").Value = "  -0.63%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.39"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.70"
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.376"
$ws.Range("E9").Value = "  -2.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0820"
$ws.Range("E10").Value = "  -1.31%  "
$ws.Range("E11").Value = "  +0.41%  "
$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.350.23"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.64"
$ws.Range("E13").Value = "  -1.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.03"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.779"
$ws.Range("E15").Value = "  +2.20%  "
$ws.Range("E16").Value = "  -1.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.062.84"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "37.848.76"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.55"
$ws.Range("E19").Value = "  -0.37%  "
$ws.Range("E20").Value = "  -4.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0823"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.26"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.43"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.26"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.12"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("E27").Value = "  +0.76%  "
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.81"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("E30").Value = "  -2.50%  "
$ws.Range("E31").Value = "  -0.44%  "
$ws.Range("E32").Value = "  +8.46%  "
$ws.Range("E33").Value = "  -2.20%  "
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0594"
$ws.Range("E35").Value = "  -2.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.60"
$ws.Range("E36").Value = "  +5.20%  "
$ws.Range("E37").Value = "  +2.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.46"
$ws.Range("E38").Value = "  +6.24%  "
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.17"
$ws.Range("E40").Value = "  +7.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.539.22"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "96.31"
$ws.Range("E43").Value = "  -1.88%  "
$ws.Range("E44").Value = "  -0.57%  "
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.19"
$ws.Range("E46").Value = "  +4.31%  "
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -0.74%  "
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.12"
$ws.Range("E50").Value = "  +1.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.238.90"
$ws.Range("E51").Value = "  -0.48%  "
